# Apply the "alternate branch" edit to the first paragraph of the document:
#   1. Append two trailing spaces to the existing sentence.
#   2. Append a new, red-colored run with the extra annotation text.
#   3. Insert a new, empty paragraph right after the (now two-run) paragraph.

$d = $word.ActiveDocument

# --- Step 1: add the two trailing spaces to the original sentence -----------
$found = $d.Content.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false,
    $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

# --- Step 2: insert the new blank paragraph right after paragraph 1 ---------
# (done before adding the colored run so the new paragraph mark does not
#  inherit the red color formatting)
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$insertionPoint = $d.Range($r1.End - 1, $r1.End - 1)
$insertionPoint.InsertParagraphAfter()

# --- Step 3: append the red "(This is a change ...)" run --------------------
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$endOfPara1 = $d.Range($r1.End - 1, $r1.End - 1)
$endOfPara1.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch alternate)")
$endOfPara1.Font.Color = 192
